{"js": "// Each entry is [oldText, newText]. The first pair updates the date\n// heading; the rest update the \"three-digit x one-digit\" answer cells\n// in the practice table.\nconst replacements = [\n  [\"2025-12-25 Thursday\", \"2025-12-26 Friday\"],\n  [\"348\u00d73=1044\", \"910\u00d78=7280\"],\n  [\"792\u00d76=4752\", \"390\u00d75=1950\"],\n  [\"133\u00d79=1197\", \"338\u00d73=1014\"],\n  [\"332\u00d77=2324\", \"349\u00d78=2792\"],\n  [\"669\u00d72=1338\", \"826\u00d74=3304\"],\n  [\"248\u00d73=744\", \"477\u00d73=1431\"],\n  [\"184\u00d74=736\", \"109\u00d78=872\"],\n  [\"785\u00d74=3140\", \"120\u00d79=1080\"],\n  [\"346\u00d72=692\", \"488\u00d75=2440\"],\n  [\"649\u00d73=1947\", \"176\u00d78=1408\"],\n  [\"205\u00d78=1640\", \"208\u00d74=832\"],\n  [\"590\u00d74=2360\", \"241\u00d77=1687\"],\n  [\"792\u00d79=7128\", \"194\u00d79=1746\"],\n  [\"887\u00d74=3548\", \"194\u00d72=388\"],\n  [\"599\u00d75=2995\", \"805\u00d72=1610\"],\n  [\"409\u00d73=1227\", \"825\u00d78=6600\"],\n  [\"601\u00d79=5409\", \"163\u00d79=1467\"],\n  [\"307\u00d74=1228\", \"570\u00d75=2850\"],\n  [\"548\u00d72=1096\", \"334\u00d79=3006\"],\n  [\"388\u00d73=1164\", \"419\u00d72=838\"],\n  [\"539\u00d78=4312\", \"497\u00d72=994\"],\n  [\"766\u00d75=3830\", \"478\u00d75=2390\"],\n  [\"451\u00d75=2255\", \"386\u00d78=3088\"],\n  [\"193\u00d74=772\", \"686\u00d75=3430\"],\n  [\"422\u00d77=2954\", \"965\u00d73=2895\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each pair is (oldText, newText). The first pair updates the date\n# heading; the rest update the \"three-digit x one-digit\" answer cells\n# in the practice table.\n$pairs = @(\n  @(\"2025-12-25 Thursday\", \"2025-12-26 Friday\"),\n  @(\"348\u00d73=1044\", \"910\u00d78=7280\"),\n  @(\"792\u00d76=4752\", \"390\u00d75=1950\"),\n  @(\"133\u00d79=1197\", \"338\u00d73=1014\"),\n  @(\"332\u00d77=2324\", \"349\u00d78=2792\"),\n  @(\"669\u00d72=1338\", \"826\u00d74=3304\"),\n  @(\"248\u00d73=744\", \"477\u00d73=1431\"),\n  @(\"184\u00d74=736\", \"109\u00d78=872\"),\n  @(\"785\u00d74=3140\", \"120\u00d79=1080\"),\n  @(\"346\u00d72=692\", \"488\u00d75=2440\"),\n  @(\"649\u00d73=1947\", \"176\u00d78=1408\"),\n  @(\"205\u00d78=1640\", \"208\u00d74=832\"),\n  @(\"590\u00d74=2360\", \"241\u00d77=1687\"),\n  @(\"792\u00d79=7128\", \"194\u00d79=1746\"),\n  @(\"887\u00d74=3548\", \"194\u00d72=388\"),\n  @(\"599\u00d75=2995\", \"805\u00d72=1610\"),\n  @(\"409\u00d73=1227\", \"825\u00d78=6600\"),\n  @(\"601\u00d79=5409\", \"163\u00d79=1467\"),\n  @(\"307\u00d74=1228\", \"570\u00d75=2850\"),\n  @(\"548\u00d72=1096\", \"334\u00d79=3006\"),\n  @(\"388\u00d73=1164\", \"419\u00d72=838\"),\n  @(\"539\u00d78=4312\", \"497\u00d72=994\"),\n  @(\"766\u00d75=3830\", \"478\u00d75=2390\"),\n  @(\"451\u00d75=2255\", \"386\u00d78=3088\"),\n  @(\"193\u00d74=772\", \"686\u00d75=3430\"),\n  @(\"422\u00d77=2954\", \"965\u00d73=2895\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $result = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $result) {\n    throw \"Replacement failed for: $old\"\n  }\n}\n"}
